$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 user settings
$ws.Range("C3").Value = "en"
$ws.Range("F3").Value = "2025-11-12 14:20:26"
$ws.Range("G3").Value = "2025-11-12 14:20:36"
$ws.Range("H3").ClearContents()
